$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: give the existing "294" group a top+bottom border (matches new border style) ---
# Pull the border/font template that already exists on row 3 (font + thin bottom border),
# then add a thin top edge so the group is fully boxed top & bottom.
$ws.Range("A3:E3").Copy()

$row6 = $ws.Range("A6:E6")
$row6.PasteSpecial(-4122)
$row6.Borders.Item(8).LineStyle = 1

# --- Row 7: brand-new data row (Gummi/Желе line) ---
$row7 = $ws.Range("A7:E7")
$row7.PasteSpecial(-4122)
$row7.Borders.Item(8).LineStyle = 1

$ws.Range("C7").Value = ' Gummi? [CS:I]Apple[CR]?\nOh, I can\''t decide…'
$ws.Range("A7").Value = "SCRIPT/P01P04A/um1105.ssb"
$ws.Range("D7").Value = ' Желе? [CS:I]Яблоко[CR]? Ох, не могу\nрешить...'
$ws.Range("E7").Value = ' Çåìå? [CS:I]Ÿáìïëï[CR]? Ïö, îå íïãô\nñåšéóû…'
$ws.Range("B7").Value = 269

$ws.Range("A7:E7").RowHeight = 43.2

$excel.CutCopyMode = $false

# --- Selection as left by the editing session ---
$ws.Range("D10").Select()
